# Weekly update: a new price-report row for "Vega Monumental Concepción -
# Berenjena" was inserted ahead of the existing series (row 156), pushing
# the previously-existing rows 156-166 down to 157-167.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 156; this shifts rows 156..166 down to 157..167
# (and carries the dimension out to R167), matching the rest of the diff.
$ws.Rows.Item(156).Insert()

# Populate the newly inserted row 156 with this week's data.
$ws.Cells.Item(156, 1).Value = 11
$ws.Cells.Item(156, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(156, 3).Value = "Bíobío"
$ws.Cells.Item(156, 4).Value = 45106
$ws.Cells.Item(156, 5).Value = 8
$ws.Cells.Item(156, 6).Value = 100112001
$ws.Cells.Item(156, 7).Value = "Berenjena"
$ws.Cells.Item(156, 8).Value = "Sin especificar"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 100
$ws.Cells.Item(156, 11).Value = 7500
$ws.Cells.Item(156, 12).Value = 8000
$ws.Cells.Item(156, 13).Value = 7750
$ws.Cells.Item(156, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(156, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(156, 16).Value = 155
$ws.Cells.Item(156, 17).Value = 50
$ws.Cells.Item(156, 18).Value = "Hortaliza"
